$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.715.48'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('D3').Value = '2.095.91'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '343.48'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.008'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5165'
$ws.Range('E7').Value = '  -1.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4381'
$ws.Range('E8').Value = '  -3.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.76'
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09283'
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.164'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.89'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').Value = '2.109.52'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.294'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.749'
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '99.48'
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001152'
$ws.Range('E17').Value = '  -1.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.009'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.78'
$ws.Range('E19').Value = '  +1.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06652'
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.198'
$ws.Range('E22').Value = '  -2.76%  '
$ws.Range('D23').Value = '29.738.86'
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.53'
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.317'
$ws.Range('E25').Value = '  -2.75%  '
$ws.Range('D26').Value = '2.345.42'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.95'
$ws.Range('E27').Value = '  -2.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.520'
$ws.Range('E28').Value = '  -4.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '161.18'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('E30').Value = '  -2.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.137'
$ws.Range('E31').Value = '  -7.11%  '
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.653'
$ws.Range('E33').Value = '  -2.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.166'
$ws.Range('E34').Value = '  -3.82%  '
$ws.Range('E35').Value = '  -2.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.268'
$ws.Range('E36').Value = '  +2.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.23'
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02577'
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06712'
$ws.Range('E39').Value = '  -3.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.49'
$ws.Range('E40').Value = '  -2.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6892'
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2222'
$ws.Range('E42').Value = '  -5.13%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.316'
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6775'
$ws.Range('E44').Value = '  +3.83%  '
$ws.Range('E45').Value = '  -3.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.321'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000356'
$ws.Range('E47').Value = '  -4.04%  '
$ws.Range('E48').Value = '  -3.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.219'
$ws.Range('E49').Value = '  -2.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.00'
$ws.Range('E50').Value = '  -2.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.160'
